$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 524 (shifts existing rows 524:627 down to 525:628,
# matching the new sheet dimension A1:R628).
$ws.Rows(524).Insert()

# Populate the newly inserted row 524 with a fresh weekly observation.
# It mirrors the row above it (row 523, "Vega Central Mapocho de Santiago"
# / Perejil / Primera) but carries a new date and its own price figures.
$ws.Range("A524").Value = 9
$ws.Range("B524").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C524").Value = "Metropolitana"
$ws.Range("D524").Value = 45209
$ws.Range("E524").Value = 13
$ws.Range("F524").Value = 100112044
$ws.Range("G524").Value = "Perejil"
$ws.Range("H524").Value = "Sin especificar"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 70
$ws.Range("K524").Value = 14000
$ws.Range("L524").Value = 15000
$ws.Range("M524").Value = 14500
$ws.Range("N524").Value = '$/docena de atados'
$ws.Range("O524").Value = "Región Metropolitana"
$ws.Range("P524").Value = 4833
$ws.Range("Q524").Value = 3
$ws.Range("R524").Value = "Hortaliza"
